# Update sample file with correct email format
# - Replace the 5 test rows (email/fullName) with 2 real rows.
# - Fix up the hyperlinks (mailto:) to match the new emails.
# - Resize columns and move the active selection, mirroring the
#   author's manual edit in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace data -------------------------------------------------
# Drop rows 4-6 entirely (only 2 data rows remain after the edit).
$ws.Range("A4:B6").Clear()

# Remove the old (stale) hyperlinks before rewriting A2:A3 so we don't
# leave dangling relationships pointing at cleared rows.
$ws.Cells.Hyperlinks.Delete()

$ws.Range("A2").Value = "minhntse140988@fpt.edu.vn"
$ws.Range("B2").Value = "Tran Nhat Minh"
$ws.Range("A3").Value = "buunqse140936@fpt.edu.vn"
$ws.Range("B3").Value = "Nguyen Quoc Buu"

# --- Hyperlinks -----------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:minhntse140988@fpt.edu.vn")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:buunqse140936@fpt.edu.vn")

# Adding a hyperlink re-applies the "Hyperlink" cell style with a brand
# new style record; force it back onto the original shared style so the
# two email cells keep looking like the rest of the workbook.
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("A3").Style = "Hyperlink"

# --- Column widths ----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 36
$ws.Columns.Item(2).ColumnWidth = 8.5

# --- Selection ----------------------------------------------------
$ws.Range("D12").Select()
